# Insert a new weekly price record row before the existing row 373
# (Hortaliza, Agrícola del Norte S.A. de Arica - Brócoli), shifting all
# subsequent rows down by one (373-419 -> 374-420).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(373).Insert()

$ws.Range("A373").Value = 1
$ws.Range("B373").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C373").Value = "Arica y Parinacota"
$ws.Range("D373").Value = 44776
$ws.Range("E373").Value = 15
$ws.Range("F373").Value = 100112023
$ws.Range("G373").Value = "Brócoli"
$ws.Range("H373").Value = "Sin especificar"
$ws.Range("I373").Value = "Tercera"
$ws.Range("J373").Value = 800
$ws.Range("K373").Value = 450
$ws.Range("L373").Value = 500
$ws.Range("M373").Value = 475
$ws.Range("N373").Value = "$/unidad"
$ws.Range("O373").Value = "Región de Arica y Parinacota"
$ws.Range("P373").Value = 475
$ws.Range("Q373").Value = 1
$ws.Range("R373").Value = "Hortaliza"
